# Update countries & provincias Spain
# Applies the 17-Jun-2020 08:33 -> 09:50 data refresh to the "Pais" sheet:
#  - reorders a few country label rows (Armenia/Japon, Groenlandia/Islas
#    Malvinas, Seychelles/Montserrat, Papua Nueva Guinea/Islas Virgenes
#    Britanicas) by swapping their label text (the underlying per-row
#    numbers stay attached to the row, exactly as the source diff shows)
#  - refreshes the numeric counters for the countries whose figures moved
#  - bumps the "Datos actualizados" timestamp footer

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country label swaps (adjacent rows trade names) ---------------------
$ws.Range("A52").Value = "Armenia"
$ws.Range("A53").Value = "Japon"

$ws.Range("A206").Value = "Groenlandia"
$ws.Range("A207").Value = "Islas Malvinas"

$ws.Range("A210").Value = "Seychelles"
$ws.Range("A211").Value = "Montserrat"

$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("A214").Value = "Islas Virgenes Britanicas"

# --- Footer timestamp (row 1) ----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 09:50"

# --- Row 4: Estados Unidos --------------------------------------------------
$ws.Range("B4").Value = 2208402
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 903042
$ws.Range("E4").Value = 1186228

# --- Row 6: Rusia ------------------------------------------------------------
$ws.Range("B6").Value = 553301
$ws.Range("C6").Value = 7843
$ws.Range("D6").Value = 304342
$ws.Range("E6").Value = 241481
$ws.Range("G6").Value = 194
$ws.Range("H6").Value = 7478

# --- Row 33: Singapur --------------------------------------------------------
$ws.Range("B33").Value = 41216
$ws.Range("C33").Value = 247
$ws.Range("E33").Value = 10027

# --- Row 52: now Armenia ------------------------------------------------------
$ws.Range("B52").Value = 18033
$ws.Range("C52").Value = 544
$ws.Range("D52").Value = 6814
$ws.Range("E52").Value = 10917
$ws.Range("G52").Value = 9
$ws.Range("H52").Value = 302

# --- Row 53: now Japon ---------------------------------------------------------
$ws.Range("B53").Value = 17587
$ws.Range("D53").Value = 15701
$ws.Range("E53").Value = 959
$ws.Range("H53").Value = 927

# --- Row 65: Chequia -------------------------------------------------------------
$ws.Range("B65").Value = 10112
$ws.Range("C65").Value = 1
$ws.Range("D65").Value = 7360
$ws.Range("E65").Value = 2420
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = 332

# --- Row 85: Hungria -----------------------------------------------------------
$ws.Range("B85").Value = 4078
$ws.Range("C85").Value = 1
$ws.Range("D85").Value = 2547
$ws.Range("E85").Value = 964
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 567

# --- Row 103: Estonia ----------------------------------------------------------
$ws.Range("B103").Value = 1977
$ws.Range("C103").Value = 2
$ws.Range("D103").Value = 1743
$ws.Range("E103").Value = 165

# --- Row 105: Mali -------------------------------------------------------------
$ws.Range("E105").Value = 634
$ws.Range("G105").Value = 2
$ws.Range("H105").Value = 106

# --- Row 125: Letonia ----------------------------------------------------------
$ws.Range("B125").Value = 1104
$ws.Range("C125").Value = 6
$ws.Range("E125").Value = 199
$ws.Range("G125").Value = 2
$ws.Range("H125").Value = 30

# --- Row 154: Taiwan -------------------------------------------------------------
$ws.Range("D154").Value = 434
$ws.Range("E154").Value = 4

# --- Row 210: now Seychelles -----------------------------------------------------
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

# --- Row 211: now Montserrat ------------------------------------------------------
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# --- Row 213: now Papua Nueva Guinea -----------------------------------------------
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

# --- Row 214: now Islas Virgenes Britanicas -----------------------------------------
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
